$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# report created: time label (shared string "11:31" -> "16:31")
# Leading apostrophe preserves the "stored as text" (quote-prefix) flag,
# matching the original cell's text formatting instead of Excel
# reinterpreting "16:31" as a time value.
$ws.Range("F24").Value = "'16:31"

# Supplier Stock 1 (J) and Supplier Unit Price 1 (P) updates per row.
# Adjusted Supplier Subtotal 1 (Q) is a formula (=P*O) and recalculates
# automatically, as do Q24 (SUM) and R24 (Q24/N25).

$ws.Range("J2").Value = 37050
$ws.Range("P2").Value = 0.019

$ws.Range("P3").Value = 0.0362

$ws.Range("J4").Value = 33872
$ws.Range("P4").Value = 0.01448

$ws.Range("J5").Value = 9999
$ws.Range("P5").Value = 0.02896

$ws.Range("J6").Value = 59043
$ws.Range("P6").Value = 0.06877

$ws.Range("J7").Value = 98272
$ws.Range("P7").Value = 0.10497

$ws.Range("J8").Value = 57880
$ws.Range("P8").Value = 0.27147

$ws.Range("J9").Value = 116092
$ws.Range("P9").Value = 0.34477

$ws.Range("J10").Value = 21003
$ws.Range("P10").Value = 0.67325

$ws.Range("J11").Value = 94764
$ws.Range("P11").Value = 0.2407

$ws.Range("J12").Value = 1836811
$ws.Range("P12").Value = 0.01176

$ws.Range("J13").Value = 68380
$ws.Range("P13").Value = 0.02896

$ws.Range("J14").Value = 83930
$ws.Range("P14").Value = 0.01176

$ws.Range("J15").Value = 131942
$ws.Range("P15").Value = 0.01176

$ws.Range("J16").Value = 98927
$ws.Range("P16").Value = 0.02896

$ws.Range("J17").Value = 891522
$ws.Range("P17").Value = 0.01448

$ws.Range("J18").Value = 16412
$ws.Range("P18").Value = 0.03529

$ws.Range("J19").Value = 14994
$ws.Range("P19").Value = 0.65244

$ws.Range("P20").Value = 3.47

$ws.Range("J21").Value = 27364
$ws.Range("P21").Value = 0.60629

$ws.Range("J22").Value = 530
$ws.Range("P22").Value = 1.23

$wb.Save()
